$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") is stored as text in this sheet (values like "1.00" or
# "0.0000171" must keep their exact literal formatting). Force the whole data
# range to a Text number format before writing, otherwise the COM input
# parser would silently coerce numeric-looking strings into real numbers and
# normalize/trim their text representation (e.g. "1.00" -> 1).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '62.063.64'
$ws.Range("E2").Value = '  +1.67%  '

# Row 3
$ws.Range("D3").Value = '3.427.99'
$ws.Range("E3").Value = '  +1.23%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = '579.00'
$ws.Range("E5").Value = '  +1.24%  '

# Row 6
$ws.Range("D6").Value = '144.85'

# Row 7
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("D9").Value = '7.59'
$ws.Range("E9").Value = '  -0.90%  '

# Row 10
$ws.Range("E10").Value = '  +1.03%  '

# Row 11
$ws.Range("D11").Value = '0.386'
$ws.Range("E11").Value = '  -0.27%  '

# Row 12
$ws.Range("D12").Value = '4.012.64'
$ws.Range("E12").Value = '  +1.21%  '

# Row 13
$ws.Range("D13").Value = '28.63'
$ws.Range("E13").Value = '  +2.46%  '

# Row 14
$ws.Range("E14").Value = '  -0.73%  '

# Row 15
$ws.Range("D15").Value = '3.420.45'

# Row 16
$ws.Range("D16").Value = '0.0000171'
$ws.Range("E16").Value = '  +0.01%  '

# Row 17
$ws.Range("D17").Value = '62.089.07'
$ws.Range("E17").Value = '  +1.57%  '

# Row 18
$ws.Range("E18").Value = '  +1.52%  '

# Row 19
$ws.Range("D19").Value = '14.05'
$ws.Range("E19").Value = '  +2.79%  '

# Row 20
$ws.Range("D20").Value = '9.23'
$ws.Range("E20").Value = '  +3.09%  '

# Row 21
$ws.Range("D21").Value = '392.74'
$ws.Range("E21").Value = '  +1.95%  '

# Row 22
$ws.Range("D22").Value = '74.77'
$ws.Range("E22").Value = '  -1.32%  '

# Row 23
$ws.Range("D23").Value = '0.555'
$ws.Range("E23").Value = '  +0.34%  '

# Row 24
$ws.Range("E24").Value = '  -0.09%  '

# Row 25
$ws.Range("D25").Value = '0.0000117'
$ws.Range("E25").Value = '  +0.84%  '

# Row 26
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '3.565.25'
$ws.Range("E26").Value = '  +1.23%  '

# Row 27
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.185'
$ws.Range("E27").Value = '  +0.21%  '

# Row 28
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '7.53'
$ws.Range("E28").Value = '  +4.21%  '

# Row 29
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.13%  '

# Row 30
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '8.03'
$ws.Range("E30").Value = '  +0.68%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '2.16'
$ws.Range("E31").Value = '  +0.85%  '

# Row 32
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '1.41'
$ws.Range("E32").Value = '  +2.75%  '

# Row 33
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.01%  '

# Row 34
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '23.62'
$ws.Range("E34").Value = '  +1.52%  '

# Row 35
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = '5.29'
$ws.Range("E35").Value = '  +5.98%  '

# Row 36
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Value = '7.00'
$ws.Range("E36").Value = '  +0.34%  '

# Row 37
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '167.73'
$ws.Range("E37").Value = '  +0.93%  '

# Row 38
$ws.Range("D38").Value = '1.51'
$ws.Range("E38").Value = '  +3.38%  '

# Row 39
$ws.Range("B39").Value = 'RenzoRestakedETH'
$ws.Range("C39").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D39").Value = '3.458.82'
$ws.Range("E39").Value = '  +1.15%  '

# Row 40
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = '28.43'
$ws.Range("E40").Value = '  +7.58%  '

# Row 41
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '0.0753'
$ws.Range("E41").Value = '  -1.99%  '

# Row 42
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = '0.789'
$ws.Range("E42").Value = '  +1.45%  '

# Row 43
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '4.45'
$ws.Range("E43").Value = '  +1.83%  '

# Row 44
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '1.68'
$ws.Range("E44").Value = '  +1.63%  '

# Row 45
$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").Value = '1.17'
$ws.Range("E45").Value = '  +4.54%  '

# Row 46
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.516.32'
$ws.Range("E46").Value = '  +2.38%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '23.04'
$ws.Range("E47").Value = '  +0.51%  '

# Row 48
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '6.67'
$ws.Range("E48").Value = '  +0.24%  '

# Row 49
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.02%  '

# Row 50
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '0.0265'
$ws.Range("E50").Value = '  +0.66%  '

# Row 51
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '2.14'
$ws.Range("E51").Value = '  +0.62%  '
